$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 81

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 44656
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100112001
$ws.Cells.Item($row, 7).Value = "Berenjena"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 220
$ws.Cells.Item($row, 11).Value = 8000
$ws.Cells.Item($row, 12).Value = 8500
$ws.Cells.Item($row, 13).Value = 8227
$ws.Cells.Item($row, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 137
$ws.Cells.Item($row, 17).Value = 60
$ws.Cells.Item($row, 18).Value = "Hortaliza"
